# Updating Assembly Manual and BOM

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 23: 18th BOM line item (BNC bulkhead connector) ---
$ws.Range("A23").Value = 18

$ws.Range("B23").Value = "2057-RF1-106-D-00-50-HDW-ND"
$ws.Range("B23").NumberFormat = "0"
$ws.Range("B23").HorizontalAlignment = -4131

$ws.Range("C23").Value = "BNC STRAIGHT BULKHEAD SKT 50 OHM"
$ws.Range("D23").Value = "DIGIKEY"
$ws.Range("E23").Value = 1

$ws.Range("F23").Value = 1.48
$ws.Range("F23").NumberFormat = "0.00"

$ws.Range("G23").Formula = "=E23*F23"
$ws.Range("G23").NumberFormat = "0.00"

$ws.Range("H23").Value = "https://www.digikey.com/en/products/detail/adam-tech/RF1-106-D-00-50-HDW/9830449"
$ws.Range("H23").NumberFormat = "0.00"

# --- Row 1: add "Total Cost:" label (D1) and total formula (E1) ---
$ws.Range("D1").Value = "Total Cost:"

$ws.Range("E1").Formula = "=SUM(G:G)"
$ws.Range("E1").NumberFormat = "0.00"

# --- Update selection to reflect the author's last active cell ---
$ws.Range("G29").Select()
